$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet Hoja1 (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.46 = 39810.79 pesos`n✅ 39810.79 pesos = 9.44 = 968.05 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 105.7
$wsTasas.Range("O10").Value = 4208

$wsTasas.Range("N12").Value = 4215.89
$wsTasas.Range("O12").Value = 102.515
